$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Forecast Comparison" sheet: the forecast data refreshed by one week.
# Each row gets a new Week_Start_Date (col B) plus refreshed forecast
# numbers (MyForecast=D, Amazon Mean=E, P70=F, P80=G, P90=H).
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

# r, Week_Start_Date, MyForecast(D), AmazonMean(E), P70(F), P80(G), P90(H)
$forecastRows = @(
    @(2,  "2025-02-02", 1, 3, 4, 5, 7),
    @(3,  "2025-02-09", 1, 3, 4, 5, 8),
    @(4,  "2025-02-16", 1, 3, 4, 5, 7),
    @(5,  "2025-02-23", 1, 3, 4, 6, 8),
    @(6,  "2025-03-02", 2, 4, 4, 6, 9),
    @(7,  "2025-03-09", 2, 4, 5, 6, 9),
    @(8,  "2025-03-16", 2, 4, 5, 7, 11),
    @(9,  "2025-03-23", 2, 4, 5, 7, 11),
    @(10, "2025-03-30", 2, 4, 5, 7, 10),
    @(11, "2025-04-06", 2, 4, 5, 7, 12),
    @(12, "2025-04-13", 2, 4, 5, 8, 12),
    @(13, "2025-04-20", 2, 5, 5, 8, 12),
    @(14, "2025-04-27", 2, 5, 5, 8, 12),
    @(15, "2025-05-04", 2, 4, 5, 7, 11),
    @(16, "2025-05-11", 2, 4, 4, 7, 11),
    @(17, "2025-05-18", 2, 4, 4, 7, 11)
)

foreach ($row in $forecastRows) {
    $r = $row[0]

    # Column B holds the date as literal text (e.g. "2025-02-02"), not a
    # real Excel date serial - force text formatting before assigning so
    # COM doesn't auto-convert the string into a date number.
    $bCell = $wsForecast.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]

    $wsForecast.Cells.Item($r, 4).Value = $row[2]   # D - MyForecast
    $wsForecast.Cells.Item($r, 5).Value = $row[3]   # E - Amazon Mean Forecast
    $wsForecast.Cells.Item($r, 6).Value = $row[4]   # F - Amazon P70 Forecast
    $wsForecast.Cells.Item($r, 7).Value = $row[5]   # G - Amazon P80 Forecast
    $wsForecast.Cells.Item($r, 8).Value = $row[6]   # H - Amazon P90 Forecast
}

# ---------------------------------------------------------------------------
# "Summary" sheet: a handful of derived metrics need to be refreshed to
# match the updated forecast window above.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

function Set-TextCell($sheet, $addr, $value) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextCell $wsSummary "B2"  "2022-12-25 to 2025-01-26"   # Historical Range
Set-TextCell $wsSummary "B9"  "23"                          # Total Forecast (16 Weeks)
Set-TextCell $wsSummary "B10" "11"                          # Total Forecast (8 Weeks)
Set-TextCell $wsSummary "B11" "5"                           # Total Forecast (4 Weeks)
Set-TextCell $wsSummary "B13" "2025-03-02"                  # Max Forecast Week
Set-TextCell $wsSummary "B15" "2025-02-02"                  # Min Forecast Week
